{"js": "// Map of old text -> new text, derived from the OOXML diff.\nconst replacements = [\n  [\"2025-09-28 Sunday\", \"2025-09-29 Monday\"],\n  [\"533\u00f73=177, 2\", \"305\u00f74=76, 1\"],\n  [\"233\u00f72=116, 1\", \"859\u00f78=107, 3\"],\n  [\"647\u00f79=71, 8\", \"315\u00f79=35, 0\"],\n  [\"168\u00f77=24, 0\", \"249\u00f74=62, 1\"],\n  [\"603\u00f77=86, 1\", \"963\u00f77=137, 4\"],\n  [\"888\u00f79=98, 6\", \"764\u00f78=95, 4\"],\n  [\"695\u00f77=99, 2\", \"566\u00f77=80, 6\"],\n  [\"192\u00f74=48, 0\", \"809\u00f72=404, 1\"],\n  [\"686\u00f77=98, 0\", \"365\u00f77=52, 1\"],\n  [\"799\u00f79=88, 7\", \"234\u00f76=39, 0\"],\n  [\"581\u00f76=96, 5\", \"242\u00f76=40, 2\"],\n  [\"780\u00f75=156, 0\", \"651\u00f77=93, 0\"],\n  [\"935\u00f76=155, 5\", \"240\u00f74=60, 0\"],\n  [\"591\u00f77=84, 3\", \"103\u00f77=14, 5\"],\n  [\"873\u00f73=291, 0\", \"463\u00f76=77, 1\"],\n  [\"877\u00f75=175, 2\", \"695\u00f75=139, 0\"],\n  [\"814\u00f77=116, 2\", \"837\u00f73=279, 0\"],\n  [\"923\u00f73=307, 2\", \"581\u00f72=290, 1\"],\n  [\"728\u00f79=80, 8\", \"340\u00f77=48, 4\"],\n  [\"833\u00f74=208, 1\", \"142\u00f76=23, 4\"],\n  [\"937\u00f74=234, 1\", \"869\u00f74=217, 1\"],\n  [\"358\u00f79=39, 7\", \"693\u00f76=115, 3\"],\n  [\"256\u00f75=51, 1\", \"678\u00f77=96, 6\"],\n  [\"164\u00f77=23, 3\", \"210\u00f75=42, 0\"],\n  [\"316\u00f75=63, 1\", \"899\u00f75=179, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Map of old text -> new text, derived from the OOXML diff.\n$replacements = @(\n    @(\"2025-09-28 Sunday\", \"2025-09-29 Monday\"),\n    @(\"533\u00f73=177, 2\", \"305\u00f74=76, 1\"),\n    @(\"233\u00f72=116, 1\", \"859\u00f78=107, 3\"),\n    @(\"647\u00f79=71, 8\", \"315\u00f79=35, 0\"),\n    @(\"168\u00f77=24, 0\", \"249\u00f74=62, 1\"),\n    @(\"603\u00f77=86, 1\", \"963\u00f77=137, 4\"),\n    @(\"888\u00f79=98, 6\", \"764\u00f78=95, 4\"),\n    @(\"695\u00f77=99, 2\", \"566\u00f77=80, 6\"),\n    @(\"192\u00f74=48, 0\", \"809\u00f72=404, 1\"),\n    @(\"686\u00f77=98, 0\", \"365\u00f77=52, 1\"),\n    @(\"799\u00f79=88, 7\", \"234\u00f76=39, 0\"),\n    @(\"581\u00f76=96, 5\", \"242\u00f76=40, 2\"),\n    @(\"780\u00f75=156, 0\", \"651\u00f77=93, 0\"),\n    @(\"935\u00f76=155, 5\", \"240\u00f74=60, 0\"),\n    @(\"591\u00f77=84, 3\", \"103\u00f77=14, 5\"),\n    @(\"873\u00f73=291, 0\", \"463\u00f76=77, 1\"),\n    @(\"877\u00f75=175, 2\", \"695\u00f75=139, 0\"),\n    @(\"814\u00f77=116, 2\", \"837\u00f73=279, 0\"),\n    @(\"923\u00f73=307, 2\", \"581\u00f72=290, 1\"),\n    @(\"728\u00f79=80, 8\", \"340\u00f77=48, 4\"),\n    @(\"833\u00f74=208, 1\", \"142\u00f76=23, 4\"),\n    @(\"937\u00f74=234, 1\", \"869\u00f74=217, 1\"),\n    @(\"358\u00f79=39, 7\", \"693\u00f76=115, 3\"),\n    @(\"256\u00f75=51, 1\", \"678\u00f77=96, 6\"),\n    @(\"164\u00f77=23, 3\", \"210\u00f75=42, 0\"),\n    @(\"316\u00f75=63, 1\", \"899\u00f75=179, 4\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
